# Apply the commit: remove the post row for "「無限の宇宙を旅した光」..."
# which was the row at sheet row 739 (A739/B739/C739). Deleting the entire
# worksheet row shifts all subsequent rows (740-815) up by one (739-814),
# matching the updated dimension ref "A1:C814".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(739).Delete()
